$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new KEYWORDS / ISQUESTION columns ---
$ws.Cells.Item(1, 8).Value = "KEYWORDS"
$ws.Cells.Item(1, 9).Value = "ISQUESTION"

# --- Row 2 (ID 1, welcome message) ---
$ws.Cells.Item(2, 2).Value = "(Welcome Message) What would you like to learn about? 1. (emoji) Exercise - brief description2. (emoji) WASH- brief description3. (emoji) Nutrition-brief description4. (emoji) Maternal Infant Care-brief description5. (emoji) Mental Health- brief description"
$ws.Cells.Item(2, 3).Value = "image"
$ws.Cells.Item(2, 4).Value = "2, 3, 4, 5, 6"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = "low data welcome"
$ws.Cells.Item(2, 8).Value = "1,2,3,4,5"
$ws.Cells.Item(2, 9).Value = $true

# --- Row 3 (ID 2, Exercise) ---
$ws.Cells.Item(3, 2).Value = "Exercise"
$ws.Cells.Item(3, 3).Value = "image 2"
$ws.Cells.Item(3, 4).Value = 7
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = "low data exercise"
$ws.Cells.Item(3, 9).Value = $false

# --- Row 4 (ID 3, WASH) ---
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "WASH"
$ws.Cells.Item(4, 3).Value = "google.com"
$ws.Cells.Item(4, 4).Value = 8
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = "low data wash"
$ws.Cells.Item(4, 9).Value = $false

# --- Row 5 (ID 4, Nutrition) ---
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Nutrition"
$ws.Cells.Item(5, 3).Value = "test1"
$ws.Cells.Item(5, 4).Value = 9
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(5, 7).Value = "low data nutrition"
$ws.Cells.Item(5, 9).Value = $false

# --- Row 6 (ID 5, Maternal Infant Care) ---
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Maternal Infant Care"
$ws.Cells.Item(6, 3).Value = "test2"
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(6, 7).Value = "low data maternal health"
$ws.Cells.Item(6, 9).Value = $false

# --- Row 7 (ID 6, Mental Health) ---
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Mental Health"
$ws.Cells.Item(7, 3).Value = "test"
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 5
$ws.Cells.Item(7, 7).Value = "low data mental halth"
$ws.Cells.Item(7, 9).Value = $false

# --- Row 8 (ID 7, next exercise) ---
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "next exercise"
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = "low data exer 2"
$ws.Cells.Item(8, 9).Value = $false

# --- Row 9 (ID 8, next wash) ---
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "next wash"
$ws.Cells.Item(9, 3).Value = "test5"
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 2
$ws.Cells.Item(9, 7).Value = "low data wash 2"
$ws.Cells.Item(9, 9).Value = $false

# --- Row 10 (ID 9, next nutrition) ---
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "next nutrition"
$ws.Cells.Item(10, 3).Value = "test7"
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(10, 7).Value = "low daya nutrition 2"
$ws.Cells.Item(10, 9).Value = $false

# --- Column widths (best achievable given 1/6-character internal rounding) ---
$ws.Columns.Item(2).ColumnWidth = 11.498697916666666
$ws.Columns.Item(3).ColumnWidth = 12.385416666666666
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).ColumnWidth = 8.608072916666666
$ws.Columns.Item(7).ColumnWidth = 12.053385416666666
$ws.Columns.Item(8).ColumnWidth = 10.721354166666666

# --- Selection / active cell ---
$ws.Range("C6").Select()
